# Applies a row-data permutation across rows 8-35 (column Artfynd data) on
# the active worksheet. Columns A, B, D, E, F, G, H, Q, R move as a unit
# from a source row to a target row; row 21 is untouched. All other
# columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY)
# stay put since they describe the observation location/date, not species.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together from the source row to the target row.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# target row -> source row (row 21 is fixed / excluded)
$mapping = @{
    8  = 27
    9  = 22
    10 = 18
    11 = 20
    12 = 25
    13 = 35
    14 = 16
    15 = 29
    16 = 10
    17 = 19
    18 = 31
    19 = 28
    20 = 30
    22 = 17
    23 = 8
    24 = 12
    25 = 32
    26 = 34
    27 = 33
    28 = 11
    29 = 24
    30 = 26
    31 = 13
    32 = 23
    33 = 9
    34 = 14
    35 = 15
}

# Snapshot the "before" values for every row that participates, so that
# overlapping source/target rows don't clobber data mid-update.
# (Use Value2 — Value's getter returns a bogus descriptor string in this
# COM shim for some property overload resolution reason; the setter is
# fine either way, but Value2 is used throughout for consistency.)
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Write the permuted values back out.
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $srcValues[$col]
    }
}
